$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.396.98'
$ws.Range("E2").Value = '  +2.10%  '
$ws.Range("D3").Value = '3.357.99'
$ws.Range("E3").Value = '  +3.65%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '''192.70'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.50%  '
$ws.Range("D6").Value = '''592.88'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.18%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  +1.16%  '
$ws.Range("D9").Value = '''0.134'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.39%  '
$ws.Range("D10").Value = '''6.73'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.74%  '
$ws.Range("D11").Value = '''0.424'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.33%  '
$ws.Range("D12").Value = '3.931.74'
$ws.Range("E12").Value = '  +3.34%  '
$ws.Range("E13").Value = '  +0.62%  '
$ws.Range("D14").Value = '''28.28'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.97%  '
$ws.Range("D15").Value = '69.439.80'
$ws.Range("E15").Value = '  +2.14%  '
$ws.Range("D16").Value = '''0.0000172'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.60%  '
$ws.Range("D17").Value = '3.328.95'
$ws.Range("E17").Value = '  +2.78%  '
$ws.Range("D18").Value = '''5.84'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.95%  '
$ws.Range("D19").Value = '''13.76'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.21%  '
$ws.Range("D20").Value = '''427.30'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +7.81%  '
$ws.Range("D21").Value = '''7.74'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.34%  '
$ws.Range("D22").Value = '''73.33'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.04%  '
$ws.Range("E23").Value = '  +0.28%  '
$ws.Range("D24").Value = '''0.518'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.80%  '
$ws.Range("D25").Value = '''0.0000122'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.25%  '
$ws.Range("E26").Value = '  +2.47%  '
$ws.Range("D27").Value = '''9.63'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.67%  '
$ws.Range("E28").Value = '  +0.06%  '
$ws.Range("E29").Value = '  +2.87%  '
$ws.Range("B30").Value = 'EthereumClassic'
$ws.Range("C30").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D30").Value = '''23.05'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.43%  '
$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D31").Value = '''5.61'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.54%  '
$ws.Range("E32").Value = '  +2.09%  '
$ws.Range("D33").Value = '''7.01'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.04%  '
$ws.Range("E34").Value = '  +0.03%  '
$ws.Range("D35").Value = '''164.95'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.77%  '
$ws.Range("D36").Value = '''1.52'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.60%  '
$ws.Range("E37").Value = '  +1.87%  '
$ws.Range("D38").Value = '''27.10'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.24%  '
$ws.Range("D39").Value = '''0.809'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.07%  '
$ws.Range("D40").Value = '''4.58'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.48%  '
$ws.Range("D41").Value = '2.754.95'
$ws.Range("E41").Value = '  +5.83%  '
$ws.Range("D42").Value = '''6.47'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.17%  '
$ws.Range("E43").Value = '  +1.99%  '
$ws.Range("D44").Value = '''25.55'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.46%  '
$ws.Range("B45").Value = 'Hedera'
$ws.Range("C45").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D45").Value = '''0.0688'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.40%  '
$ws.Range("B46").Value = 'OKB'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D46").Value = '''41.12'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.17%  '
$ws.Range("D47").Value = '''344.18'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.45%  '
$ws.Range("E48").Value = '  +1.41%  '
$ws.Range("D49").Value = '''32.63'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.01%  '
$ws.Range("D50").Value = '''1.01'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.34%  '
$ws.Range("D51").Value = '''6.29'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.05%  '
